$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "43.699.12"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.288.44"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'102.57"
$ws.Range("E5").Value = "  +5.50%  "
$ws.Range("D6").Value = "'270.44"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("D10").Value = "'46.10"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("D11").Value = "'0.0934"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'8.05"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "'15.54"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'0.855"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "2.285.91"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "43.672.02"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'72.24"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'2.50"
$ws.Range("E21").Value = "  +10.66%  "
$ws.Range("D22").Value = "'233.23"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("D23").Value = "'2.88"
$ws.Range("E23").Value = "  +13.34%  "
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").Value = "'41.66"
$ws.Range("E27").Value = "  +8.88%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'177.41"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").Value = "'21.80"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'5.47"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  +10.19%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").Value = "'0.0358"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "'3.55"
$ws.Range("E38").Value = "  +4.79%  "
$ws.Range("D39").Value = "'0.236"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").Value = "'1.38"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'65.69"
$ws.Range("E42").Value = "  +4.88%  "
$ws.Range("D43").Value = "'12.22"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'5.26"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "'8.79"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").Value = "'98.86"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  +6.73%  "
$ws.Range("D50").Value = "'1.53"
$ws.Range("E50").Value = "  +10.90%  "
$ws.Range("D51").Value = "2.511.12"
$ws.Range("E51").Value = "  -1.27%  "
